$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 329
$ws.Range("F3").Value = 1144
$ws.Range("F4").Value = 1244
$ws.Range("F5").Value = 1126
$ws.Range("F6").Value = 3396
$ws.Range("F7").Value = 63
$ws.Range("F8").Value = 1184
$ws.Range("F10").Value = 599
$ws.Range("F12").Value = 157
$ws.Range("F13").Value = 654
$ws.Range("F14").Value = 1821
$ws.Range("F15").Value = 55
$ws.Range("F16").Value = 405
$ws.Range("F17").Value = 58
$ws.Range("F18").Value = 73
$ws.Range("F19").Value = 688
$ws.Range("F20").Value = 456
$ws.Range("F22").Value = 805
$ws.Range("F23").Value = 80113
$ws.Range("F24").Value = 80113
$ws.Range("F26").Value = 672
$ws.Range("F27").Value = 33878
$ws.Range("F28").Value = 33878
$ws.Range("F29").Value = 538
$ws.Range("F30").Value = 30
$ws.Range("F32").Value = 59
$ws.Range("F33").Value = 52
$ws.Range("F34").Value = 1004
$ws.Range("F35").Value = 313
$ws.Range("F37").Value = 631
$ws.Range("F38").Value = 2680
$ws.Range("F39").Value = 2680
$ws.Range("F40").Value = 1217
$ws.Range("F41").Value = 5503
$ws.Range("F42").Value = 796
$ws.Range("F43").Value = 459
$ws.Range("F47").Value = 426
$ws.Range("F50").Value = 19
$ws.Range("F51").Value = 54
$ws.Range("F52").Value = 7

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 9
$ws.Range("F11").Value = 1982
$ws.Range("F12").Value = 31
$ws.Range("F14").Value = 83
$ws.Range("F15").Value = 418
$ws.Range("F16").Value = 11
$ws.Range("F17").Value = 76
$ws.Range("F20").Value = 535
$ws.Range("F41").Value = 36
$ws.Range("F44").Value = 73
$ws.Range("G44").Value = 224
$ws.Range("F45").Value = 828
$ws.Range("F46").Value = 208
$ws.Range("F48").Value = 70

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 105
$ws.Range("F4").Value = 735
$ws.Range("F5").Value = 583
$ws.Range("F6").Value = 612
$ws.Range("F7").Value = 159

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 329
$ws.Range("F4").Value = 1144
$ws.Range("F5").Value = 1126
$ws.Range("F6").Value = 3397
$ws.Range("F7").Value = 1184
$ws.Range("F9").Value = 612
$ws.Range("F10").Value = 612
$ws.Range("F13").Value = 599
$ws.Range("F15").Value = 157
$ws.Range("F16").Value = 654
$ws.Range("F17").Value = 159
$ws.Range("F18").Value = 1821
$ws.Range("F19").Value = 31
$ws.Range("F20").Value = 55
$ws.Range("F21").Value = 405
$ws.Range("F22").Value = 58
$ws.Range("F23").Value = 73
$ws.Range("F24").Value = 805
$ws.Range("F25").Value = 11
$ws.Range("F26").Value = 76
$ws.Range("F27").Value = 80115
$ws.Range("F28").Value = 672
$ws.Range("F29").Value = 33879
$ws.Range("F30").Value = 538
$ws.Range("F31").Value = 30
$ws.Range("F33").Value = 535
$ws.Range("F34").Value = 535
$ws.Range("F35").Value = 52
$ws.Range("F38").Value = 313
$ws.Range("F41").Value = 2680
$ws.Range("F42").Value = 1217
$ws.Range("F43").Value = 5503
$ws.Range("F44").Value = 796
$ws.Range("F49").Value = 426
$ws.Range("F51").Value = 73
$ws.Range("G51").Value = 224
$ws.Range("F53").Value = 208
$ws.Range("F54").Value = 54
$ws.Range("F55").Value = 7
